# Update countries & provincias Spain
# - Refresh the "last updated" timestamp
# - Update Australia's case counts
# - Finlandia overtakes Arabia Saudita / Indonesia in the ranking (new Finlandia
#   numbers are inserted right after Sudafrica, pushing Arabia Saudita and
#   Indonesia down one row each)
# - Bulgaria overtakes Ucrania in the ranking (new Bulgaria numbers are inserted
#   right after Hungria, pushing Ucrania down one row)
# - Update Banglades' active/recovered counts

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados..." timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 28 de Marzo de 2020 a las 08:59"

# Australia (row 21) - updated case numbers
$ws.Range("B21").Value = 3635
$ws.Range("C21").Value = 257
$ws.Range("E21").Value = 3451

# Finlandia jumps ahead of Arabia Saudita and Indonesia (rows 38-40)
$ws.Range("A38").Value = "Finlandia"
$ws.Range("B38").Value = 1163
$ws.Range("C38").Value = 122
$ws.Range("D38").Value = 10
$ws.Range("E38").Value = 1146
$ws.Range("F38").Value = 32
$ws.Range("G38").Value = 0
$ws.Range("H38").Value = 7

$ws.Range("A39").Value = "Arabia Saudita"
$ws.Range("B39").Value = 1104
$ws.Range("C39").Value = 0
$ws.Range("D39").Value = 35
$ws.Range("E39").Value = 1066
$ws.Range("F39").Value = 6
$ws.Range("G39").Value = 0
$ws.Range("H39").Value = 3

$ws.Range("A40").Value = "Indonesia"
$ws.Range("B40").Value = 1046
$ws.Range("C40").Value = 0
$ws.Range("D40").Value = 46
$ws.Range("E40").Value = 913
$ws.Range("F40").Value = 0
$ws.Range("G40").Value = 0
$ws.Range("H40").Value = 87

# Bulgaria jumps ahead of Ucrania (rows 71-72)
$ws.Range("A71").Value = "Bulgaria"
$ws.Range("B71").Value = 313
$ws.Range("C71").Value = 20
$ws.Range("D71").Value = 9
$ws.Range("E71").Value = 300
$ws.Range("F71").Value = 8
$ws.Range("G71").Value = 1
$ws.Range("H71").Value = 4

$ws.Range("A72").Value = "Ucrania"
$ws.Range("B72").Value = 310
$ws.Range("C72").Value = 0
$ws.Range("D72").Value = 5
$ws.Range("E72").Value = 300
$ws.Range("F72").Value = 0
$ws.Range("G72").Value = 0
$ws.Range("H72").Value = 5

# Banglades (row 125) - updated active/recovered counts
$ws.Range("D125").Value = 15
$ws.Range("E125").Value = 28
